# Applies the coin price/volume/hour updates described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'291.29"
$ws.Range("E2").Value = "'-6.79%"
$ws.Range("G2").Value = "'4"
$ws.Range("D3").Value = "'40.50"
$ws.Range("E3").Value = "'1.27%"
$ws.Range("G3").Value = "'4"
$ws.Range("D4").Value = "'5.009"
$ws.Range("E4").Value = "'-2.01%"
$ws.Range("G4").Value = "'4"
$ws.Range("D5").Value = "'0.07329"
$ws.Range("E5").Value = "'-3.24%"
$ws.Range("G5").Value = "'4"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.531"
$ws.Range("E6").Value = "'-8.69%"
$ws.Range("G6").Value = "'4"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9225"
$ws.Range("E7").Value = "'-0.77%"
$ws.Range("G7").Value = "'4"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.363"
$ws.Range("E8").Value = "'-2.52%"
$ws.Range("G8").Value = "'4"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1211"
$ws.Range("E9").Value = "'0.12%"
$ws.Range("G9").Value = "'4"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1731"
$ws.Range("E10").Value = "'-4.50%"
$ws.Range("G10").Value = "'4"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08663"
$ws.Range("E11").Value = "'-4.72%"
$ws.Range("G11").Value = "'4"
$ws.Range("D12").Value = "'0.04291"
$ws.Range("E12").Value = "'3.52%"
$ws.Range("G12").Value = "'4"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("E13").Value = "'0.14%"
$ws.Range("G13").Value = "'4"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001267"
$ws.Range("E14").Value = "'-0.91%"
$ws.Range("G14").Value = "'4"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005803"
$ws.Range("E15").Value = "'-0.60%"
$ws.Range("G15").Value = "'4"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.336"
$ws.Range("E16").Value = "'-0.45%"
$ws.Range("G16").Value = "'4"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.295"
$ws.Range("E17").Value = "'-0.34%"
$ws.Range("G17").Value = "'4"
$ws.Range("E18").Value = "'-1.97%"
$ws.Range("G18").Value = "'4"
$ws.Range("D19").Value = "'7.721"
$ws.Range("E19").Value = "'1.06%"
$ws.Range("G19").Value = "'4"
$ws.Range("E20").Value = "'2.93%"
$ws.Range("G20").Value = "'4"
$ws.Range("D21").Value = "'0.2796"
$ws.Range("E21").Value = "'-0.50%"
$ws.Range("G21").Value = "'4"
$ws.Range("D22").Value = "'0.03939"
$ws.Range("E22").Value = "'-1.77%"
$ws.Range("G22").Value = "'4"
$ws.Range("E23").Value = "'-0.39%"
$ws.Range("G23").Value = "'4"
$ws.Range("D24").Value = "'0.003773"
$ws.Range("E24").Value = "'-6.81%"
$ws.Range("G24").Value = "'4"
$ws.Range("E25").Value = "'0.93%"
$ws.Range("G25").Value = "'4"
$ws.Range("D26").Value = "'0.0003729"
$ws.Range("G26").Value = "'4"
$ws.Range("G27").Value = "'4"
$ws.Range("G28").Value = "'4"
$ws.Range("G29").Value = "'4"
$ws.Range("G30").Value = "'4"
$ws.Range("G31").Value = "'4"
$ws.Range("G32").Value = "'4"
$ws.Range("G33").Value = "'4"
$ws.Range("G34").Value = "'4"
$ws.Range("G35").Value = "'4"
$ws.Range("G36").Value = "'4"
$ws.Range("G37").Value = "'4"
$ws.Range("D38").Value = "'0.02284"
$ws.Range("E38").Value = "'-6.08%"
$ws.Range("G38").Value = "'4"
$ws.Range("D39").Value = "'0.04974"
$ws.Range("E39").Value = "'-3.61%"
$ws.Range("G39").Value = "'4"
$ws.Range("B40").Value = "KickToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D40").Value = "'0.007682"
$ws.Range("E40").Value = "'-0.20%"
$ws.Range("G40").Value = "'4"
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D41").Value = "'0.005153"
$ws.Range("E41").Value = "'135.61%"
$ws.Range("G41").Value = "'4"
$ws.Range("E42").Value = "'-1.27%"
$ws.Range("G42").Value = "'4"
$ws.Range("D43").Value = "'0.007377"
$ws.Range("E43").Value = "'-3.08%"
$ws.Range("G43").Value = "'4"
$ws.Range("D44").Value = "'0.007928"
$ws.Range("E44").Value = "'-2.58%"
$ws.Range("G44").Value = "'4"
$ws.Range("D45").Value = "'0.3158"
$ws.Range("E45").Value = "'1.44%"
$ws.Range("G45").Value = "'4"
$ws.Range("D46").Value = "'0.00006375"
$ws.Range("E46").Value = "'-3.30%"
$ws.Range("G46").Value = "'4"
$ws.Range("E47").Value = "'0.16%"
$ws.Range("G47").Value = "'4"
$ws.Range("D48").Value = "'0.02038"
$ws.Range("E48").Value = "'-92.43%"
$ws.Range("G48").Value = "'4"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("E49").Value = "'0.16%"
$ws.Range("G49").Value = "'4"
$ws.Range("D50").Value = "'0.0002004"
$ws.Range("E50").Value = "'0.16%"
$ws.Range("G50").Value = "'4"
$ws.Range("G51").Value = "'4"
